$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original (unstyled) cell formatting; setting a leading
# apostrophe forces Excel to store these numeric-looking price strings
# as literal text (preventing silent numeric conversion / precision loss),
# but it also stamps the cells with a "quote prefix" style, so we capture
# the original style here and re-apply it to the whole price column below.
$origDStyle = $ws.Range("D2").Style

# --- Price (column D) updates ---
$ws.Range("D2").Value = "'22.399.15"
$ws.Range("D3").Value = "'1.572.51"
$ws.Range("D6").Value = "'290.47"
$ws.Range("D7").Value = "'0.3764"
$ws.Range("D8").Value = "'49.84"
$ws.Range("D9").Value = "'0.3427"
$ws.Range("D10").Value = "'0.07648"
$ws.Range("D14").Value = "'6.024"
$ws.Range("D15").Value = "'6.945"
$ws.Range("D16").Value = "'1.570.37"
$ws.Range("D18").Value = "'90.28"
$ws.Range("D19").Value = "'0.06761"
$ws.Range("D21").Value = "'16.81"
$ws.Range("D22").Value = "'6.203"
$ws.Range("D24").Value = "'22.402.30"
$ws.Range("D25").Value = "'2.412"
$ws.Range("D26").Value = "'2.686"
$ws.Range("D28").Value = "'147.48"
$ws.Range("D29").Value = "'5.035"
$ws.Range("D30").Value = "'126.29"
$ws.Range("D31").Value = "'1.746.02"
$ws.Range("D32").Value = "'6.169"
$ws.Range("D33").Value = "'2.008"
$ws.Range("D34").Value = "'0.9960"
$ws.Range("D35").Value = "'9.991"
$ws.Range("D36").Value = "'0.08573"
$ws.Range("D39").Value = "'0.06577"
$ws.Range("D41").Value = "'5.442"
$ws.Range("D43").Value = "'0.6417"
$ws.Range("D45").Value = "'14.06"
$ws.Range("D46").Value = "'3.787"
$ws.Range("D47").Value = "'0.5994"
$ws.Range("D49").Value = "'2.088"
$ws.Range("D51").Value = "'0.07328"

# Restore original styling for the price column
$ws.Range("D2:D51").Style = $origDStyle

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -10.43%  "
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  +5.88%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +7.78%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  +0.46%  "

